$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three score cells in row 24 from 2 to 5
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 5

# Update the selected/active cell to G24
$ws.Range("G24").Select()
